# Append 6 new daily rows (27 Apr 2021 - 02 May 2021) to the data table,
# matching the style of the existing rows (date in column A styled with
# the same date format/border/bold as the rows above, plain numbers in
# columns B/C/D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti
$rows = @(
    @(44313, 8, 57, 172.8162993057029),
    @(44314, 1, 53, 160.6888397053027),
    @(44315, 5, 54, 163.7207046054028),
    @(44316, 8, 44, 133.4020556044023),
    @(44317, 5, 43, 130.3701907043022),
    @(44318, 5, 39, 118.242731103902)
)

$startRow = 239
$lastExistingRow = $startRow - 1

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Copy the formatting of the last existing date cell (column A) down
    # onto the new row so the new cell picks up the same style index
    # (bold, bordered, centered, date-formatted).
    $ws.Cells.Item($lastExistingRow, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}

$excel.CutCopyMode = $false
